# "cols wrong in template"
# The helper block of cells (Test/RISE/TIME metadata used for image placement)
# in a few rows of the template was shifted one column to the left of where
# it belongs (columns M:O instead of N:P like the rest of the sheet). This
# script shifts those cells one column to the right, and resets the saved
# scroll position of the sheet view back to the top (it had drifted to A5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Shift-ColumnsRight {
    param(
        $ws,
        [int]$rowStart,
        [int]$rowEnd,
        [int]$colStart,
        [int]$colEnd
    )

    # Snapshot values + horizontal alignment (the only formatting this
    # particular block of cells uses) before touching anything, since the
    # source and destination ranges overlap (shift by a single column).
    $vals = @{}
    $aligns = @{}
    for ($r = $rowStart; $r -le $rowEnd; $r++) {
        for ($c = $colStart; $c -le $colEnd; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $vals["$r,$c"] = $cell.Value2
            $aligns["$r,$c"] = $cell.HorizontalAlignment
        }
    }

    # Wipe the whole source block (contents + formatting) now that it has
    # been captured.
    $srcRange = $ws.Range($ws.Cells.Item($rowStart, $colStart), $ws.Cells.Item($rowEnd, $colEnd))
    $srcRange.Clear()

    # Re-create each originally-populated cell one column to the right.
    for ($r = $rowStart; $r -le $rowEnd; $r++) {
        for ($c = $colStart; $c -le $colEnd; $c++) {
            $v = $vals["$r,$c"]
            $al = $aligns["$r,$c"]
            $hadValue = ($v -ne $null -and $v -ne "")
            $hadStyle = ($al -ne 1)   # 1 == xlGeneral (the unstyled default)
            if ($hadValue -or $hadStyle) {
                $destCell = $ws.Cells.Item($r, $c + 1)
                if ($hadValue) {
                    $destCell.Value2 = $v
                }
                if ($hadStyle) {
                    $destCell.HorizontalAlignment = $al
                }
            }
        }
    }
}

# Self Test section (rows 66-70): helper cells were in M:O, belong in N:P.
Shift-ColumnsRight $ws 66 70 13 15

# Sample Rate and Delay Time Accuracy section (rows 101-102): same mistake.
Shift-ColumnsRight $ws 101 102 13 15

# The saved view had scrolled to row 5 (topLeftCell="A5"); restore it to the
# top of the sheet.
$excel.ActiveWindow.ScrollRow = 1
